$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Post Treatment" values for column D (rows 2-12)
$values = @(64, 67, 67, 54, 57, 56, 64, 63, 60, 63, 63)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Update the selected cell to D13, matching the saved selection state
$ws.Range("D13").Select()
